$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_G2 = New-Object 'object[,]' 1,4
$arr_G2[0,0] = 19.539082
$arr_G2[0,1] = 58.61724600000001
$arr_G2[0,2] = 0.224220971665117
$arr_G2[0,3] = 0.224220971665117
$ws.Range("G2:J2").Value = $arr_G2
$arr_M2 = New-Object 'object[,]' 1,8
$arr_M2[0,0] = 1.174933333333333
$arr_M2[0,1] = 3.5248
$arr_M2[0,2] = 0.01171850713626266
$arr_M2[0,3] = 0.01171850713626266
$arr_M2[0,4] = 22.95711874453334
$arr_M2[0,5] = 206.6140687008
$arr_M2[0,6] = 0.00262753505655742
$arr_M2[0,7] = 0.00262753505655742
$ws.Range("M2:T2").Value = $arr_M2

$arr_G3 = New-Object 'object[,]' 1,4
$arr_G3[0,0] = 19.539082
$arr_G3[0,1] = 58.61724600000001
$arr_G3[0,2] = 0.224220971665117
$arr_G3[0,3] = 0.224220971665117
$ws.Range("G3:J3").Value = $arr_G3
$arr_O3 = New-Object 'object[,]' 1,6
$arr_O3[0,0] = 0.2743256641287217
$arr_O3[0,1] = 0.2743256641287218
$arr_O3[0,2] = 537.4171618318054
$arr_O3[0,3] = 4836.754456486249
$arr_O3[0,4] = 0.0615095669636205
$arr_O3[0,5] = 0.06150956696362052
$ws.Range("O3:T3").Value = $arr_O3

$arr_G4 = New-Object 'object[,]' 1,4
$arr_G4[0,0] = 19.539082
$arr_G4[0,1] = 58.61724600000001
$arr_G4[0,2] = 0.224220971665117
$arr_G4[0,3] = 0.224220971665117
$ws.Range("G4:J4").Value = $arr_G4
$arr_M4 = New-Object 'object[,]' 1,8
$arr_M4[0,0] = 39.361408
$arr_M4[0,1] = 118.084224
$arr_M4[0,2] = 0.3925813724534833
$arr_M4[0,3] = 0.3925813724534833
$arr_M4[0,4] = 769.0857785474562
$arr_M4[0,5] = 6921.772006927105
$arr_M4[0,6] = 0.08802497678914521
$arr_M4[0,7] = 0.08802497678914521
$ws.Range("M4:T4").Value = $arr_M4

$arr_G5 = New-Object 'object[,]' 1,4
$arr_G5[0,0] = 19.539082
$arr_G5[0,1] = 58.61724600000001
$arr_G5[0,2] = 0.224220971665117
$arr_G5[0,3] = 0.224220971665117
$ws.Range("G5:J5").Value = $arr_G5
$arr_M5 = New-Object 'object[,]' 1,8
$arr_M5[0,0] = 32.221985
$arr_M5[0,1] = 96.665955
$arr_M5[0,2] = 0.3213744562815322
$arr_M5[0,3] = 0.3213744562815322
$arr_M5[0,4] = 629.5880071177701
$arr_M5[0,5] = 5666.29206405993
$arr_M5[0,6] = 0.0720588928557938
$arr_M5[0,7] = 0.0720588928557938
$ws.Range("M5:T5").Value = $arr_M5

$arr_I6 = New-Object 'object[,]' 1,2
$arr_I6[0,0] = 0.3010605798326856
$arr_I6[0,1] = 0.3010605798326856
$ws.Range("I6:J6").Value = $arr_I6
$arr_M6 = New-Object 'object[,]' 1,8
$arr_M6[0,0] = 1.174933333333333
$arr_M6[0,1] = 3.5248
$arr_M6[0,2] = 0.01171850713626266
$arr_M6[0,3] = 0.01171850713626266
$arr_M6[0,4] = 30.82442926364444
$arr_M6[0,5] = 277.4198633728
$arr_M6[0,6] = 0.0035279805532167
$arr_M6[0,7] = 0.003527980553216699
$ws.Range("M6:T6").Value = $arr_M6

$arr_I7 = New-Object 'object[,]' 1,2
$arr_I7[0,0] = 0.3010605798326856
$arr_I7[0,1] = 0.3010605798326856
$ws.Range("I7:J7").Value = $arr_I7
$arr_O7 = New-Object 'object[,]' 1,2
$arr_O7[0,0] = 0.2743256641287217
$arr_O7[0,1] = 0.2743256641287218
$ws.Range("O7:P7").Value = $arr_O7
$arr_S7 = New-Object 'object[,]' 1,2
$arr_S7[0,0] = 0.08258864350557954
$arr_S7[0,1] = 0.08258864350557955
$ws.Range("S7:T7").Value = $arr_S7

$arr_I8 = New-Object 'object[,]' 1,2
$arr_I8[0,0] = 0.3010605798326856
$arr_I8[0,1] = 0.3010605798326856
$ws.Range("I8:J8").Value = $arr_I8
$arr_M8 = New-Object 'object[,]' 1,8
$arr_M8[0,0] = 39.361408
$arr_M8[0,1] = 118.084224
$arr_M8[0,2] = 0.3925813724534833
$arr_M8[0,3] = 0.3925813724534833
$arr_M8[0,4] = 1032.648323263829
$arr_M8[0,5] = 9293.834909374464
$arr_M8[0,6] = 0.1181907756223572
$arr_M8[0,7] = 0.1181907756223572
$ws.Range("M8:T8").Value = $arr_M8

$arr_I9 = New-Object 'object[,]' 1,2
$arr_I9[0,0] = 0.3010605798326856
$arr_I9[0,1] = 0.3010605798326856
$ws.Range("I9:J9").Value = $arr_I9
$arr_M9 = New-Object 'object[,]' 1,8
$arr_M9[0,0] = 32.221985
$arr_M9[0,1] = 96.665955
$arr_M9[0,2] = 0.3213744562815322
$arr_M9[0,3] = 0.3213744562815322
$arr_M9[0,4] = 845.3452372049866
$arr_M9[0,5] = 7608.10713484488
$arr_M9[0,6] = 0.09675318015153216
$arr_M9[0,7] = 0.09675318015153217
$ws.Range("M9:T9").Value = $arr_M9

$arr_G10 = New-Object 'object[,]' 1,4
$arr_G10[0,0] = 19.67155566666667
$arr_G10[0,1] = 59.014667
$arr_G10[0,2] = 0.2257411748281949
$arr_G10[0,3] = 0.2257411748281949
$ws.Range("G10:J10").Value = $arr_G10
$arr_M10 = New-Object 'object[,]' 1,8
$arr_M10[0,0] = 1.174933333333333
$arr_M10[0,1] = 3.5248
$arr_M10[0,2] = 0.01171850713626266
$arr_M10[0,3] = 0.01171850713626266
$arr_M10[0,4] = 23.11276647128889
$arr_M10[0,5] = 208.0148982416
$arr_M10[0,6] = 0.002645349568172518
$arr_M10[0,7] = 0.002645349568172518
$ws.Range("M10:T10").Value = $arr_M10

$arr_G11 = New-Object 'object[,]' 1,4
$arr_G11[0,0] = 19.67155566666667
$arr_G11[0,1] = 59.014667
$arr_G11[0,2] = 0.2257411748281949
$arr_G11[0,3] = 0.2257411748281949
$ws.Range("G11:J11").Value = $arr_G11
$arr_O11 = New-Object 'object[,]' 1,6
$arr_O11[0,0] = 0.2743256641287217
$arr_O11[0,1] = 0.2743256641287218
$arr_O11[0,2] = 541.0608141772663
$arr_O11[0,3] = 4869.547327595395
$arr_O11[0,4] = 0.06192659770594245
$arr_O11[0,5] = 0.06192659770594246
$ws.Range("O11:T11").Value = $arr_O11

$arr_G12 = New-Object 'object[,]' 1,4
$arr_G12[0,0] = 19.67155566666667
$arr_G12[0,1] = 59.014667
$arr_G12[0,2] = 0.2257411748281949
$arr_G12[0,3] = 0.2257411748281949
$ws.Range("G12:J12").Value = $arr_G12
$arr_M12 = New-Object 'object[,]' 1,8
$arr_M12[0,0] = 39.361408
$arr_M12[0,1] = 118.084224
$arr_M12[0,2] = 0.3925813724534833
$arr_M12[0,3] = 0.3925813724534833
$arr_M12[0,4] = 774.3001285903788
$arr_M12[0,5] = 6968.701157313408
$arr_M12[0,6] = 0.08862178023331449
$arr_M12[0,7] = 0.08862178023331449
$ws.Range("M12:T12").Value = $arr_M12

$arr_G13 = New-Object 'object[,]' 1,4
$arr_G13[0,0] = 19.67155566666667
$arr_G13[0,1] = 59.014667
$arr_G13[0,2] = 0.2257411748281949
$arr_G13[0,3] = 0.2257411748281949
$ws.Range("G13:J13").Value = $arr_G13
$arr_M13 = New-Object 'object[,]' 1,8
$arr_M13[0,0] = 32.221985
$arr_M13[0,1] = 96.665955
$arr_M13[0,2] = 0.3213744562815322
$arr_M13[0,3] = 0.3213744562815322
$arr_M13[0,4] = 633.8565716179983
$arr_M13[0,5] = 5704.709144561984
$arr_M13[0,6] = 0.07254744732076544
$arr_M13[0,7] = 0.07254744732076544
$ws.Range("M13:T13").Value = $arr_M13

$arr_G14 = New-Object 'object[,]' 1,4
$arr_G14[0,0] = 21.69639766666667
$arr_G14[0,1] = 65.089193
$arr_G14[0,2] = 0.2489772736740025
$arr_G14[0,3] = 0.2489772736740025
$ws.Range("G14:J14").Value = $arr_G14
$arr_M14 = New-Object 'object[,]' 1,8
$arr_M14[0,0] = 1.174933333333333
$arr_M14[0,1] = 3.5248
$arr_M14[0,2] = 0.01171850713626266
$arr_M14[0,3] = 0.01171850713626266
$arr_M14[0,4] = 25.49182083182222
$arr_M14[0,5] = 229.4263874864
$arr_M14[0,6] = 0.002917641958316018
$arr_M14[0,7] = 0.002917641958316018
$ws.Range("M14:T14").Value = $arr_M14

$arr_G15 = New-Object 'object[,]' 1,4
$arr_G15[0,0] = 21.69639766666667
$arr_G15[0,1] = 65.089193
$arr_G15[0,2] = 0.2489772736740025
$arr_G15[0,3] = 0.2489772736740025
$ws.Range("G15:J15").Value = $arr_G15
$arr_O15 = New-Object 'object[,]' 1,6
$arr_O15[0,0] = 0.2743256641287217
$arr_O15[0,1] = 0.2743256641287218
$arr_O15[0,2] = 596.7535453300316
$arr_O15[0,3] = 5370.781907970284
$arr_O15[0,4] = 0.06830085595357922
$arr_O15[0,5] = 0.06830085595357924
$ws.Range("O15:T15").Value = $arr_O15

$arr_G16 = New-Object 'object[,]' 1,4
$arr_G16[0,0] = 21.69639766666667
$arr_G16[0,1] = 65.089193
$arr_G16[0,2] = 0.2489772736740025
$arr_G16[0,3] = 0.2489772736740025
$ws.Range("G16:J16").Value = $arr_G16
$arr_M16 = New-Object 'object[,]' 1,8
$arr_M16[0,0] = 39.361408
$arr_M16[0,1] = 118.084224
$arr_M16[0,2] = 0.3925813724534833
$arr_M16[0,3] = 0.3925813724534833
$arr_M16[0,4] = 854.0007606879147
$arr_M16[0,5] = 7686.006846191231
$arr_M16[0,6] = 0.09774383980866641
$arr_M16[0,7] = 0.09774383980866641
$ws.Range("M16:T16").Value = $arr_M16

$arr_G17 = New-Object 'object[,]' 1,4
$arr_G17[0,0] = 21.69639766666667
$arr_G17[0,1] = 65.089193
$arr_G17[0,2] = 0.2489772736740025
$arr_G17[0,3] = 0.2489772736740025
$ws.Range("G17:J17").Value = $arr_G17
$arr_M17 = New-Object 'object[,]' 1,8
$arr_M17[0,0] = 32.221985
$arr_M17[0,1] = 96.665955
$arr_M17[0,2] = 0.3213744562815322
$arr_M17[0,3] = 0.3213744562815322
$arr_M17[0,4] = 699.1010001693683
$arr_M17[0,5] = 6291.909001524315
$arr_M17[0,6] = 0.08001493595344078
$arr_M17[0,7] = 0.0800149359534408
$ws.Range("M17:T17").Value = $arr_M17
